# Insert a new blank row above row 1, shifting the existing threshold
# table (rows 1-27) down to rows 2-28, then restore the active-cell
# selection to A7 (matches the post-edit sheetView/selection in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

$ws.Range("A7").Select()
